# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values computed for rows 2-21 (column G), replacing the prior Strike# values
$kValues = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 2
    6  = 2
    7  = 3
    8  = 2
    9  = 0
    10 = 2
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 4
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
